# Applies the "Cleans data and updated naming" edit to meta_data_names sheet:
#  - Removes the "Catch_units" / "Units used for survey" row content
#  - Inserts a new "Fit_0no_1yes" description row (F7:G7)
#  - Shifts the remaining F/G rows (old rows 8-14) down into (new rows 8-14),
#    with the "weight1_Numbers2" row relocating from old row 14 to new row 10
#  - Preserves the per-row cell formatting that travels with the moved content
#  - Updates the sheet view (scroll position / selection)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- 1. Fix up cell formatting (styles) before the content changes below ----
# Old F14 (style s="3") moves to F10.
$ws.Range("F14").Copy()
$ws.Range("F10").PasteSpecial(-4122)

# Old F13 (style s="2") moves to F14.
$ws.Range("F13").Copy()
$ws.Range("F14").PasteSpecial(-4122)

# Old F13 becomes a normal/unstyled cell again (copy format from neighboring
# unstyled cell F12).
$ws.Range("F12").Copy()
$ws.Range("F13").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ---- 2. Update the cell content for rows 7 through 14 in columns F and G ----
$ws.Range("F7").Value  = "Fit_0no_1yes"
$ws.Range("G7").Value  = "Index of wether data should be included in the likelihood and associated parameters estimated."

$ws.Range("F8").Value  = "Selectivity"
$ws.Range("G8").Value  = "Selectivity to use for the species: 0 = empirical selectivity provided in srv_emp_sel; 1 = logistic selectivity; 2 = non-parametric selecitivty sensu Ianelli et al 2018; 3 = double logistic"

$ws.Range("F9").Value  = "Nselages"
$ws.Range("G9").Value  = "Number of ages to estimate non-parametric selectivity for Selectivity = 2. Not used otherwise"

$ws.Range("F10").Value = "weight1_Numbers2"
$ws.Range("G10").Value = "Is the observation in weight (kg) set as 1, if the observation is in numbers caught, set as 2"

$ws.Range("F11").Value = "Weight_index"
$ws.Range("G11").Value = "Weight-at-age (wt) index to use for calculation of derived quantities"

$ws.Range("F12").Value = "ALK_index"
$ws.Range("G12").Value = "Age transition matrix (e.g. Age Length Key or ALK) index to use for derived quantitied"

$ws.Range("F13").Value = "Estimate_q"
$ws.Range("G13").Value = "Estimate catchability? (0 = no; 1 = yes, 2 = analytical from Ludwig and Walters 1994)"

$ws.Range("F14").Value = "log_q_start"
$ws.Range("G14").Value = "Starting value or fixed value for catchability"

# ---- 3. Update the sheet view: scroll back to the top and move the selection ----
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("G14").Select()
